# Change module name for MAP
# 1) Refresh the cached date placeholder text (27/01/2022 -> 22/03/2024)
#    on the slide master, every slide layout and the notes master.
# 2) Rename "AllineaForge" to "ArmForge" in the command-line textbox on
#    the "Startup" slide, and shrink the textbox to fit the new text.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "27/01/2022") {
                $tr.Text = "22/03/2024"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Rename AllineaForge -> ArmForge on the "Startup" slide and resize the
# textbox that holds it.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf("AllineaForge")
            if ($idx -ge 0) {
                $chars = $tr.Characters($idx + 1, 12)
                $chars.Text = "ArmForge"
                $shp.Width = 204.37653543307084
            }
        }
    }
}
